$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-10 from 45221 to 45224
foreach ($row in 2..10) {
    $ws.Cells.Item($row, 3).Value = 45224
}
